# Update the "want-to-go" counts (column F) for the first two events.
# These values live in both the "展览" sheet and the "全部类型" sheet,
# which mirror the same source rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 475
    $ws.Range("F3").Value = 53
}
